$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.13904857635498
$ws.Range("B1").Value = 3.329093456268311
$ws.Range("C1").Value = 3.51341438293457
$ws.Range("D1").Value = 3.940408945083618
$ws.Range("E1").Value = 1.159700751304626
